# Auto-generated script to apply Sargatanas_Profits.xlsx market-data refresh
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 482.6905
$ws.Range("J33").Value = 2949
$ws.Range("L33").Value = 2949
$ws.Range("N33").Value = -3407
$ws.Range("H55").Value = 183.66667
$ws.Range("I55").Value = 149.14285
$ws.Range("K55").Value = 149.14285
$ws.Range("M55").Value = 64.85714999999999
$ws.Range("H76").Value = 4009670.8
$ws.Range("J76").Value = 9112.5
$ws.Range("L76").Value = 9112.5
$ws.Range("N76").Value = -9742.5
$ws.Range("H79").Value = 4009670.8
$ws.Range("J79").Value = 9112.5
$ws.Range("L79").Value = 9112.5
$ws.Range("N79").Value = -11296.5
$ws.Range("H86").Value = 211112930
$ws.Range("I86").Value = 333335260
$ws.Range("K86").Value = 333335260
$ws.Range("M86").Value = -333334137
$ws.Range("H89").Value = 211112930
$ws.Range("I89").Value = 333335260
$ws.Range("K89").Value = 1666676300
$ws.Range("M89").Value = -1666670684
$ws.Range("H96").Value = 1025.1666
$ws.Range("I96").Value = 830.2
$ws.Range("K96").Value = 2490.6
$ws.Range("M96").Value = -1117.6
$ws.Range("H107").Value = 18481154
$ws.Range("J107").Value = 18184508
$ws.Range("L107").Value = 18184508
$ws.Range("N107").Value = -18188348
$ws.Range("H112").Value = 13454.5625
$ws.Range("J112").Value = 13454.5625
$ws.Range("L112").Value = 40363.6875
$ws.Range("N112").Value = -42579.6875
$ws.Range("H137").Value = 2984.775
$ws.Range("I137").Value = 3000.5454
$ws.Range("K137").Value = 9001.636200000001
$ws.Range("M137").Value = -6451.636200000001
$ws.Range("H138").Value = 1927615.4
$ws.Range("I138").Value = 4044.2856
$ws.Range("K138").Value = 12132.8568
$ws.Range("M138").Value = -6992.856800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2708.182
$ws.Range("I2").Value = 2363.077
$ws.Range("K2").Value = 2363.077
$ws.Range("M2").Value = -2250.077
$ws.Range("H32").Value = 1842253
$ws.Range("I32").Value = 2276304.5
$ws.Range("J32").Value = 5881.846
$ws.Range("K32").Value = 2276304.5
$ws.Range("L32").Value = 5881.846
$ws.Range("M32").Value = -2276017.5
$ws.Range("N32").Value = -6455.846
$ws.Range("H74").Value = 30367.027
$ws.Range("J74").Value = 6997.4546
$ws.Range("L74").Value = 6997.4546
$ws.Range("N74").Value = -8745.454600000001
$ws.Range("H77").Value = 30367.027
$ws.Range("J77").Value = 6997.4546
$ws.Range("L77").Value = 34987.273
$ws.Range("N77").Value = -43723.273
$ws.Range("H102").Value = 2128.7334
$ws.Range("I102").Value = 2128.7334
$ws.Range("K102").Value = 2128.7334
$ws.Range("M102").Value = -506.7334000000001
$ws.Range("H110").Value = 30304440
$ws.Range("I110").Value = 1295.75
$ws.Range("K110").Value = 1295.75
$ws.Range("M110").Value = 749.25
$ws.Range("H116").Value = 2708.182
$ws.Range("I116").Value = 2363.077
$ws.Range("K116").Value = 2363.077
$ws.Range("M116").Value = -69.07700000000023

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2708.182
$ws.Range("I3").Value = 2363.077
$ws.Range("K3").Value = 2363.077
$ws.Range("M3").Value = -2249.077
$ws.Range("H97").Value = 3351.7778
$ws.Range("I97").Value = 2461.875
$ws.Range("K97").Value = 2461.875
$ws.Range("M97").Value = -1470.875
$ws.Range("H134").Value = 6001.9062
$ws.Range("I134").Value = 1096.7142
$ws.Range("K134").Value = 3290.1426
$ws.Range("M134").Value = -755.1425999999997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 903
$ws.Range("I94").Value = 1971.8334
$ws.Range("J94").Value = 597.619
$ws.Range("K94").Value = 1971.8334
$ws.Range("L94").Value = 597.619
$ws.Range("M94").Value = -1520.8334
$ws.Range("N94").Value = -1499.619
$ws.Range("H134").Value = 5211.2856
$ws.Range("I134").Value = 2831.6206
$ws.Range("K134").Value = 8494.861800000001
$ws.Range("M134").Value = -5959.861800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 28266732
$ws.Range("I4").Value = 17230074
$ws.Range("K4").Value = 51690222
$ws.Range("M4").Value = -51690110
$ws.Range("H17").Value = 2984.5
$ws.Range("I17").Value = 594
$ws.Range("J17").Value = 5375
$ws.Range("K17").Value = 1782
$ws.Range("L17").Value = 16125
$ws.Range("M17").Value = -1613
$ws.Range("N17").Value = -16463
$ws.Range("H33").Value = 215.92857
$ws.Range("I33").Value = 73.333336
$ws.Range("J33").Value = 322.875
$ws.Range("K33").Value = 440.000016
$ws.Range("L33").Value = 1937.25
$ws.Range("M33").Value = -157.000016
$ws.Range("N33").Value = -2503.25
$ws.Range("H86").Value = 299
$ws.Range("I86").Value = 299
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 897
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 289
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 299
$ws.Range("I89").Value = 299
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 2691
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 3237
$ws.Range("N89").ClearContents()
$ws.Range("H107").Value = 8335481
$ws.Range("I107").Value = 674.75
$ws.Range("K107").Value = 2024.25
$ws.Range("M107").Value = -104.25
$ws.Range("H131").Value = 3579.818
$ws.Range("I131").Value = 3000
$ws.Range("J131").Value = 3637.8
$ws.Range("K131").Value = 9000
$ws.Range("L131").Value = 10913.4
$ws.Range("M131").Value = -3960
$ws.Range("N131").Value = -20993.4
$ws.Range("H140").Value = 252196.31
$ws.Range("I140").Value = 365112.9
$ws.Range("K140").Value = 1095338.7
$ws.Range("M140").Value = -1090158.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 48818.75
$ws.Range("J135").Value = 48818.75
$ws.Range("L135").Value = 48818.75
$ws.Range("N135").Value = -58958.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2085897.8
$ws.Range("I61").Value = 3572674
$ws.Range("J61").Value = 4410.95
$ws.Range("K61").Value = 3572674
$ws.Range("L61").Value = 4410.95
$ws.Range("M61").Value = -3572472
$ws.Range("N61").Value = -4814.95
$ws.Range("H113").Value = 2085897.8
$ws.Range("I113").Value = 3572674
$ws.Range("J113").Value = 4410.95
$ws.Range("K113").Value = 3572674
$ws.Range("L113").Value = 4410.95
$ws.Range("M113").Value = -3570504
$ws.Range("N113").Value = -8750.950000000001
$ws.Range("H139").Value = 77393
$ws.Range("J139").Value = 77393
$ws.Range("L139").Value = 77393
$ws.Range("N139").Value = -87673

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2233.3333
$ws.Range("I96").Value = 2233.3333
$ws.Range("K96").Value = 2233.3333
$ws.Range("M96").Value = -860.3332999999998
$ws.Range("H122").Value = 5575.1665
$ws.Range("J122").Value = 5399
$ws.Range("L122").Value = 16197
$ws.Range("N122").Value = -21097
$ws.Range("H132").Value = 29431672
$ws.Range("I132").Value = 50004652
$ws.Range("J132").Value = 41698
$ws.Range("K132").Value = 150013956
$ws.Range("L132").Value = 125094
$ws.Range("M132").Value = -150011426
$ws.Range("N132").Value = -130154
